$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before row 630 (2026/12/29 ... currently there).
# Copy the row above (629, same date block "2026/01/13") so the new row's
# date/weekday cells keep their original text formatting (not reinterpreted
# as a date serial), then insert it and overwrite the "time" value.
$ws.Range("A629:D629").Copy()
$ws.Range("A630").EntireRow.Insert()
$ws.Range("A630").PasteSpecial()

$ws.Range("C630").Value2 = 23

Write-Output "Row inserted at 630; new dimension should be A1:D672"
